$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

$ws.Cells.Item($row, 4).Value = 44656
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino dulce"
$ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 15000
$ws.Cells.Item($row, 12).Value = 16000
$ws.Cells.Item($row, 13).Value = 15500
$ws.Cells.Item($row, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 861
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
